$d = $word.ActiveDocument

# Locate the paragraph that ends the document with the old (pre-edit) wording.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*the shape of data will be lost.*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$rng = $target.Range
$rng.MoveEnd(1, 1) | Out-Null

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman" w:hint="eastAsia"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>F</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">rom the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>histograms shown, we would like to see my segmenting according to activities work. The Acceleration and angular velocity are fitted in normal distribution</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">. According to the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">box plots of activity = 1 and activity = 0, there are many outliers analyzed by a </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>traditional</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> method</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> to check the outlier</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">. In this </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>case,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> the dynamic method</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>such as rolling window and smoothing have more statistics meaning in time series.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> If we just replace all the outliers using </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">traditional </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve">methods, the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>pattern and structure</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> of data will be lost.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> The preprocessing method I used is to keep the true signal and characteristics of the data</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t xml:space="preserve"> in the time series.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>\</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a5"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:ind w:firstLineChars="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Product Overview</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a5"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>APP1:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The API is used for programmers to access </w:t>
      </w:r>
      <w:r>
        <w:t>the dataset which</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> use</w:t>
      </w:r>
      <w:r>
        <w:t>d</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the </w:t>
      </w:r>
      <w:r>
        <w:t>IMU</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:r>
        <w:t>inertial measurement unit)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to collect </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Acceleration and </w:t>
      </w:r>
      <w:r>
        <w:t>A</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>ngular</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> velocity</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> over </w:t>
      </w:r>
      <w:r>
        <w:t>a period</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Developers could integrate the </w:t>
      </w:r>
      <w:r>
        <w:t>API</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> into their own software. The API will allow developers to request and manipulate the dataset effectively. </w:t>
      </w:r>
      <w:r>
        <w:t>The product will concentrate on improving performance and allow users to update in real time. The API can play an important role in many applications</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> like mapping software and </w:t>
      </w:r>
      <w:r>
        <w:t>health monitoring</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> software on the mobile phone</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> because it could </w:t>
      </w:r>
      <w:r>
        <w:t>provide the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">data basis to </w:t>
      </w:r>
      <w:r>
        <w:t>developers</w:t>
      </w:r>
      <w:r>
        <w:t>. The developers can</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> predict the </w:t>
      </w:r>
      <w:r>
        <w:t>speed, routing,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and location </w:t>
      </w:r>
      <w:r>
        <w:t>of</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> user </w:t>
      </w:r>
      <w:r>
        <w:t>over a certain time</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/>
        </w:rPr>
        <w:t>Ap</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>p2:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/>
        </w:rPr>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">he App2 will deploy a machine learning model according to the acceleration and rotation speed in the dataset </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">by App1 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">to predict </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>whether</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> a person is moving. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">he model will also predict the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">speed, rotation angel </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">and </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">moving </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">routing of a person if they are </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>walking.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>To</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> give a detailed description </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>of the routing, the moving trajectory</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> of a person moving will be visualized on </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">a </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>3-D s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>patial</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>map.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> The predicted </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>results (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>speed and rotation angle)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> will </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>be exported</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> as csv also.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>The App2 will concentrate on how to</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> demonstrate these data to users with visual impact and</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> design a </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>user-friendly</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>interface</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> especially for somebody who does not have prior programming knowledge</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> and physics knowledge</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> The App2 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>aims to</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> assist</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>individual athletes or sports teams</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> to evaluate their performance when they are</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> moving. The trainer can collect athletes’ acceleration</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">angular speed </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">by IMU </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">to draw player performance statistics </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">through the App2 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>and help athletes to train or perform better</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> in the competition</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Times New Roman" w:hint="eastAsia"/>
          <w:color w:val="000000" w:themeColor="text1"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
Write-Host "Last paragraph text:" $d.Paragraphs.Last.Range.Text
